$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 80, shifting existing rows 80-93 down to 81-94.
$ws.Rows("80:80").Insert()

# Populate the newly inserted row 80 with the new data record.
$ws.Range("A80").Value = 11
$ws.Range("B80").Value = "Vega Monumental Concepción"
$ws.Range("C80").Value = "Bíobío"
$ws.Range("D80").Value = 44474
$ws.Range("E80").Value = 8
$ws.Range("F80").Value = 100112003
$ws.Range("G80").Value = "Ajo"
$ws.Range("H80").Value = "Chino"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 400
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 16000
$ws.Range("M80").Value = 15500
$ws.Range("N80").Value = "$/caja 10 kilos"
$ws.Range("O80").Value = "China"
$ws.Range("P80").Value = 1550
$ws.Range("Q80").Value = 10
$ws.Range("R80").Value = "Hortaliza"
